$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-51 (Coin, Link, Price, Volume(1h)).
# A new row (OKB) was inserted after the existing "Cardano" row, which
# shifted every following coin down by one row; the last row
# (previously "Decentraland") falls off the bottom of the 50-row table.
# D and E columns hold numeric-looking text (e.g. "26.051.94", "0.9980",
# "  +0.89%  ") that must stay as literal text, so values are written
# with a leading apostrophe to force the text type (preserving trailing
# zeros and surrounding spaces), and the cell style is reset to "Normal"
# afterward so no stray quote-prefix formatting is left behind.
$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '26.051.94', '  +0.89%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.747.81', '  +0.42%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '0.9988', '  -0.20%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '234.61', '  +3.58%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.9993', '  -0.10%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.5289', '  +2.46%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2781', '  +2.20%  '),
    @(9, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '40.01', '  +3.25%  '),
    @(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06179', '  +1.56%  '),
    @(11, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07253', '  +3.54%  '),
    @(12, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.744.23', '  +0.14%  '),
    @(13, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '15.35', '  +1.38%  '),
    @(14, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6442', '  +2.29%  '),
    @(15, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.608', '  +2.58%  '),
    @(16, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '78.41', '  +2.84%  '),
    @(17, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.9991', '  -0.11%  '),
    @(18, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '0.9980', '  -0.26%  '),
    @(19, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '25.981.11', '  +0.53%  '),
    @(20, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '11.63', '  +1.73%  '),
    @(21, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000006750', '  +2.21%  '),
    @(22, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '1.972.76', '  +0.56%  '),
    @(23, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.318', '  +6.41%  '),
    @(24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.823', '  +4.94%  '),
    @(25, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '5.221', '  +2.76%  '),
    @(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '139.17', '  +1.76%  '),
    @(27, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.512', '  +0.58%  '),
    @(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '15.34', '  +2.62%  '),
    @(29, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.813', '  -0.13%  '),
    @(30, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '104.51', '  +1.87%  '),
    @(31, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.08315', '  +0.11%  '),
    @(32, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '3.802', '  +5.10%  '),
    @(33, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.678', '  +9.06%  '),
    @(34, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04537', '  +3.17%  '),
    @(35, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.641', '  +1.26%  '),
    @(36, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.003', '  +3.76%  '),
    @(37, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.6323', '  +6.25%  '),
    @(38, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.708', '  +1.17%  '),
    @(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01600', '  +3.03%  '),
    @(40, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.940', '  +0.34%  '),
    @(41, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9984', '  -0.11%  '),
    @(42, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '98.38', '  -2.91%  '),
    @(43, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.3916', '  +3.11%  '),
    @(44, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.7393', '  +2.11%  '),
    @(45, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '5.059', '  +3.83%  '),
    @(46, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1143', '  +4.02%  '),
    @(47, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '6.347', '  +2.69%  '),
    @(48, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05343', '  -2.54%  '),
    @(49, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '54.05', '  +4.16%  '),
    @(50, 'Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '30.69', '  +3.13%  '),
    @(51, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '7.684', '  +3.72%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $coin = $row[1]
    $link = $row[2]
    $price = $row[3]
    $volume = $row[4]

    $ws.Cells.Item($r, 2).Value = $coin
    $ws.Cells.Item($r, 3).Value = $link

    $ws.Cells.Item($r, 4).Value = "'" + $price
    $ws.Cells.Item($r, 4).Style = "Normal"

    $ws.Cells.Item($r, 5).Value = "'" + $volume
    $ws.Cells.Item($r, 5).Style = "Normal"
}
